# Apply "a lot of figures" edit to the offloading workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("small_t0")

# --- Block 1: rows 1-20, columns A:I ---
# C: 250 -> 100, E: 400 -> 200, H: 250 -> 200
for ($r = 1; $r -le 20; $r++) {
    $ws.Cells.Item($r, 3).Value = 100   # C
    $ws.Cells.Item($r, 5).Value = 200   # E
    $ws.Cells.Item($r, 8).Value = 200   # H
}

# --- Block 2: rows 22-41, columns A:M ---
# C: 250 -> 100, E: 400 -> 200, F: 100 -> 400, G: 300 -> 400,
# J: 0.64 -> 4, K: 32 -> 44, L: 250 -> 200
for ($r = 22; $r -le 41; $r++) {
    $ws.Cells.Item($r, 3).Value = 100    # C
    $ws.Cells.Item($r, 5).Value = 200    # E
    $ws.Cells.Item($r, 6).Value = 400    # F
    $ws.Cells.Item($r, 7).Value = 400    # G
    $ws.Cells.Item($r, 10).Value = 4     # J
    $ws.Cells.Item($r, 11).Value = 44    # K
    $ws.Cells.Item($r, 12).Value = 200   # L
}

# --- Block 3: rows 43-62, columns A:O ---
# C: 250 -> 100, E: 400 -> 200, F: 100 -> 400, G: 300 -> 400,
# K: 0.35 -> 3, L: 17.5 -> 24, M: 35 -> 45, N: 250 -> 200
for ($r = 43; $r -le 62; $r++) {
    $ws.Cells.Item($r, 3).Value = 100    # C
    $ws.Cells.Item($r, 5).Value = 200    # E
    $ws.Cells.Item($r, 6).Value = 400    # F
    $ws.Cells.Item($r, 7).Value = 400    # G
    $ws.Cells.Item($r, 11).Value = 3     # K
    $ws.Cells.Item($r, 12).Value = 24    # L
    $ws.Cells.Item($r, 13).Value = 45    # M
    $ws.Cells.Item($r, 14).Value = 200   # N
}

# --- View / selection update ---
# Move the selection to L15 (was M14).
$ws.Range("L15").Select()
